$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvoiceNonPO")
$ws.Select()
